# Apply the four textual edits described by the diff.
$d = $word.ActiveDocument

# 1. "Longhedge, Salisbury," was split across two runs (with proofErr
#    spell-check markers around "Longhedge"); merge into a single run of
#    text.
$d.Content.Find.Execute("Longhedge, Salisbury,", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Longhedge, Salisbury,", 2)

# 2. " (Grade Pending)" becomes " D*D*D" (grade result replaces the
#    "pending" placeholder). Scope the search to the BTEC Nationals
#    qualification line specifically - "(Grade Pending)" also appears
#    (unchanged) after "GCSE English" elsewhere in the document.
$d.Content.Find.Execute("Nationals in Information Technology (Grade Pending)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Nationals in Information Technology D*D*D", 2)

# 3. ", Oct 2023 - Present" becomes ", Oct 2023 – Aug 2024" (role ended).
$d.Content.Find.Execute(", Oct 2023 - Present", $true, $false, $false, $false, $false,
                         $true, 1, $false, ", Oct 2023 – Aug 2024", 2)

# 4. "Assisting Customers and Manual Labour Activities" was split across
#    three runs (with proofErr spell-check markers around "Labour"); merge
#    into a single run of text.
$d.Content.Find.Execute("Assisting Customers and Manual Labour Activities", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Assisting Customers and Manual Labour Activities", 2)
